# Apply the commit changes:
# - Added "AS5048A check" style new small table (rows 8-9: Motor_RPM / Motor_Steps_per_Revo / Microsteps_Period)
# - New RPM formula converting RPM to Microstepping period
# - Updated B3 input value from 0.137 to 1
# - Updated current selection to G4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Motor_RPM input value (B3) used by the existing formulas.
$ws.Range("B3").Value = 1

# New header row (row 8) reusing existing shared strings plus one new string.
$ws.Range("B8").Value = "Motor_RPM"
$ws.Range("C8").Value = "Motor_Steps_per_Revo"
$ws.Range("D8").Value = "Microsteps_Period"

# New data row (row 9) with the RPM -> Microstepping period formula.
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 51200
$ws.Range("D9").Formula = "=(B9*C9)/(60*0.715)"

# Update the active selection to match the authored state.
$ws.Range("G4").Select()
